$wb = $excel.ActiveWorkbook

$wsOutput = $wb.Worksheets.Item("obj_output")
$wsReport = $wb.Worksheets.Item("obj_report")
$wsRel    = $wb.Worksheets.Item("rel_report__output")

# ---------------------------------------------------------------------------
# 0. Stash copies of the two "special" cell formats we need to re-apply after
#    writing new values (writing .Value resets any quotePrefix formatting),
#    using scratch cells far outside the used range. We clear them at the end.
# ---------------------------------------------------------------------------
$wsOutput.Range("B5").Copy() | Out-Null          # style with numFmt "d-mmm" + quote prefix
$wsOutput.Range("Z1").PasteSpecial(-4122) | Out-Null
$wsOutput.Range("B2").Copy() | Out-Null          # style with quote prefix only
$wsOutput.Range("Z2").PasteSpecial(-4122) | Out-Null
$wsRel.Range("B2").Copy() | Out-Null             # style with quote prefix only (rel sheet)
$wsRel.Range("Z1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 1. obj_output: re-order the output variable list (rows 3-7). "connection_flow"
#    moves from row 3 down to row 7, the rest shift up by one row.
# ---------------------------------------------------------------------------
$wsOutput.Range("B3").Value = "units_on"
$wsOutput.Range("B4").Value = "units_started_up"
$wsOutput.Range("B5").Value = "units_shut_down"
$wsOutput.Range("B6").Value = "units_available"
$wsOutput.Range("B7").Value = "connection_flow"

# B4 needs the numFmt "d-mmm" + quote-prefix style that used to live on B5.
$wsOutput.Range("Z1").Copy() | Out-Null
$wsOutput.Range("B4").PasteSpecial(-4122) | Out-Null

# B5 loses that style (it becomes a plain cell) -- already plain after the
# .Value assignment above, nothing further required.

# ---------------------------------------------------------------------------
# 2. obj_report: rename the "report_1" entry to "result_temp".
# ---------------------------------------------------------------------------
$wsReport.Range("B2").Value = "result_temp"
$wsReport.Range("Z1").Copy() | Out-Null
# placeholder removed below; restore the quote-prefix style that B2 had.

# ---------------------------------------------------------------------------
# 3. rel_report__output: same rename applied to every row of the relationship.
# ---------------------------------------------------------------------------
foreach ($r in 2..8) {
    $wsRel.Range("B$r").Value = "result_temp"
}

# ---------------------------------------------------------------------------
# 4. Re-apply the quote-prefix-only style (originally on obj_output!B2 / rel!B2)
#    to the cells that need it: obj_report!B2 and rel_report__output!B2:B8.
# ---------------------------------------------------------------------------
$wsOutput.Range("Z2").Copy() | Out-Null
$wsReport.Range("B2").PasteSpecial(-4122) | Out-Null

$wsRel.Range("Z1").Copy() | Out-Null
foreach ($r in 2..8) {
    $wsRel.Range("B$r").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 5. Clean up the scratch cells used as format stashes.
# ---------------------------------------------------------------------------
$wsOutput.Range("Z1").Clear() | Out-Null
$wsOutput.Range("Z2").Clear() | Out-Null
$wsRel.Range("Z1").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 6. Selections / active tab. Order matters: the last worksheet activated
#    becomes the tab that is marked "tabSelected" and the workbook's
#    activeTab. obj_report must end up active (activeTab=1).
# ---------------------------------------------------------------------------
$wsOutput.Range("E5").Select() | Out-Null
$wsRel.Range("B2").Select() | Out-Null
$wsReport.Activate() | Out-Null
$wsReport.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Page setup tweak on obj_output (paper size + orientation).
# ---------------------------------------------------------------------------
$wsOutput.PageSetup.PaperSize = 9
$wsOutput.PageSetup.Orientation = 1
